# Adds excel preview generation to frontend:
#  - Expand the "Simple Property Survey" sheet from a 3-column table
#    (Address / Size (SF) / Asking Rate) into a full 8-column survey table
#    with a numbered row column, Divisibility, NNN Asking Rate, Opex,
#    Direct/Sublease and Comments columns, styled header row, print area
#    and a page footer.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# column width correction factor used by this runtime: ColumnWidth (character
# units) is stored with an extra 0.8333333333333333 (5/6) of "padding" baked
# in on top of what we ask for, so subtract it to land on an exact width.
$cw = 0.8333333333333333

# ---------------------------------------------------------------------
# Column widths / row height
# ---------------------------------------------------------------------
$ws.Columns("A").ColumnWidth = 3 - $cw
$ws.Columns("B").ColumnWidth = 30 - $cw
$ws.Range("C1:G1").EntireColumn.ColumnWidth = 14 - $cw
$ws.Columns("H").ColumnWidth = 50 - $cw
$ws.Rows(1).RowHeight = 30

# ---------------------------------------------------------------------
# Header row (row 1)
# ---------------------------------------------------------------------
$ws.Range("A1").Value = ""

$headerRange = $ws.Range("B1:H1")
$headerRange.Font.Bold = $true
$headerRange.Font.ThemeColor = 2
$headerRange.Font.Size = 12
$headerRange.Interior.Color = 11954948
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108
$headerRange.WrapText = $true

$ws.Range("B1").Value = "Address"
$ws.Range("C1").Value = "Size (SF)"
$ws.Range("D1").Value = "Divisibility (SF)"
$ws.Range("E1").Value = "NNN Asking Rate (SF/Mo)"
$ws.Range("F1").Value = "Opex (SF/Mo)"
$ws.Range("G1").Value = "Direct/Sublease"
$ws.Range("H1").Value = "Comments"

# ---------------------------------------------------------------------
# Styling helpers for the data rows (rows 2-4)
# ---------------------------------------------------------------------
$rowNumRange = $ws.Range("A2:A4")
$rowNumRange.Interior.Color = 11954948
$rowNumRange.HorizontalAlignment = -4108
$rowNumRange.VerticalAlignment = -4108

$dataRange = $ws.Range("B2:G4")
$dataRange.HorizontalAlignment = -4108
$dataRange.VerticalAlignment = -4108
$dataRange.WrapText = $true

$commentsRange = $ws.Range("H2:H4")
$commentsRange.WrapText = $true

# force all data cells to be stored as text, matching the source workbook
$ws.Range("A2:H4").NumberFormat = "@"

# ---------------------------------------------------------------------
# Row 2 - 401 Lambert Avenue
# ---------------------------------------------------------------------
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "401 Lambert Avenue, Palo Alto, CA 94306"
$ws.Range("C2").Value = "8000"
$ws.Range("D2").Value = "3500 - 3500"
$ws.Range("E2").Value = "$4.00"
$ws.Range("F2").Value = ""
$ws.Range("G2").Value = "Direct Lease"
$ws.Range("H2").Value = " - Available within 60 days`n - Open Floor Plan w/ 1 conference room (can build more)`n - Kitchenette`n - Two restrooms and a Shower`n - 12 Parking Spaces in secure private garage Private outdoor balconies`n - Call for pricing"

# ---------------------------------------------------------------------
# Row 3 - 4101 El Camino Way
# ---------------------------------------------------------------------
$ws.Range("A3").Value = "2"
$ws.Range("B3").Value = "4101 El Camino Way, Palo Alto, CA 94306"
$ws.Range("C3").Value = "8975"
$ws.Range("D3").Value = "2768 - 2768"
$ws.Range("E3").Value = "$4.50"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = "Direct Lease"
$ws.Range("H3").Value = " - Rare stand-along retail building on El Camino Real`n - 11 On-site parking spaces`n - Restaurant infrastructure in place (fume hood, multiple walk-in refrigeration units)`n - Tenant Improvements are available"

# ---------------------------------------------------------------------
# Row 4 - 366 Cambridge Ave
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "3"
$ws.Range("B4").Value = "366 Cambridge Ave, Palo Alto, CA 94306"
$ws.Range("C4").Value = "4029"
$ws.Range("D4").Value = "702 - 2717"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = ""
$ws.Range("G4").Value = "Direct Lease"
$ws.Range("H4").Value = " - After Hours HVAC Available, Air Conditioning, Balcony, Bicycle Storage, CCTV (Closed Circuit`n - Television Monitoring), Central Heating, Conference Rooms, Hardwood Floors, Kitchen, Natural Light,`n - Plug & Play, Private Restrooms, Security System, Wi-Fi"

# ---------------------------------------------------------------------
# Print area + footer
# ---------------------------------------------------------------------
$ws.PageSetup.PrintArea = '$B2:$H20'
$ws.PageSetup.LeftFooter = "Page &P of &N"
